# kicad: Fix resistor value on remote
#
# Adds the missing LED (D1) and its current-limiting resistor (R2) to the
# remote's bill of materials (Sheet1), and refreshes the blank placeholder
# areas on Sheet2/Sheet3 so their formatting matches the rest of the
# workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1: append two new BOM rows (D1 LED, R2 resistor) ---------------

# Row 16 is the last populated BOM row; use it as a formatting template so
# the new rows pick up the same cell styles (fills/borders) as the rest of
# the table.
$ws1.Range("A16:I16").Copy($ws1.Range("A17:I17"))
$ws1.Range("A16:I16").Copy($ws1.Range("A18:I18"))

# Row 16's "F" (Description/Value) cell happens to be blank in the
# template, so its format doesn't carry the same style as the populated
# text cells. Pull that look from a neighboring populated cell (D16) before
# writing the new description values so the new rows render consistently.
$ws1.Range("D16").Copy()
$ws1.Range("F17").PasteSpecial(-4122)
$ws1.Range("D16").Copy()
$ws1.Range("F18").PasteSpecial(-4122)

# Row 17: D1 LED
$ws1.Range("A17").Value = "11"
$ws1.Range("B17").Value = "D1"
$ws1.Range("C17").Value = 1
$ws1.Range("D17").Value = ""
$ws1.Range("E17").Value = ""
$ws1.Range("F17").Value = "LED (any color)"
$ws1.Range("G17").Value = "0603"
$ws1.Range("H17").Value = "SMD"
$ws1.Range("I17").Value = ""

# Row 18: R2 330 Ohm resistor
$ws1.Range("A18").Value = "12"
$ws1.Range("B18").Value = "R2"
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = ""
$ws1.Range("E18").Value = ""
$ws1.Range("F18").Value = "330 Ohm"
$ws1.Range("G18").Value = "0603"
$ws1.Range("H18").Value = "SMD"
$ws1.Range("I18").Value = ""

# --- Sheet2 / Sheet3: restyle the blank A1:E10 placeholder area ----------
# These cells previously used a borderless style; give them the same
# filled/bordered look used elsewhere in the workbook (e.g. Sheet1!D9).

$ws1.Range("D9").Copy()
$ws2.Range("A1:E10").PasteSpecial(-4122)

$ws1.Range("D9").Copy()
$ws3.Range("A1:E10").PasteSpecial(-4122)
